$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix B10 to use MAX instead of SUM, matching the "Max" row label and C10's formula
$ws.Range("B10").Formula = "=MAX(B2:B6)"

# Update the selected/active cell in the sheet view to C10
$ws.Range("C10").Select()
